# edit.ps1 - applies three content changes to the active document:
#   1. "EL CUAL TIENEN LAS SIGUIENTES MEDIDA" -> "EL CUAL TIENE LAS SIGUIENTES MEDIDA"
#   2. Insert a new empty paragraph (Prrafodelista style) after the paragraph
#      ending "...EL MEDIO DE CONSIGNACION." and before "EL PAGO DEBERA REALIZARSE..."
#   3. " OBLIGADOS" -> " OBLIGAD{{SEXO_11}}S" (split across runs, wrapped with
#      proofErr gramStart/gramEnd) in the "...QUEDARAN OBLIGADOS A CUBRIR LA PENA
#      CONVENCIONAL." sentence.
#
# Whole paragraphs are rebuilt and pushed back in with Range.InsertXML so the
# exact run layout (including the new w:proofErr markers) matches precisely,
# rather than relying on Find/Replace (which silently coalesces adjacent runs
# that share identical formatting).

$d = $word.ActiveDocument

function Set-ParagraphXml($anchorText, $newParagraphXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "anchor text not found: $anchorText"
    }
    $p = $rng.Paragraphs(1)
    $prng = $p.Range
    $pkg = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$newParagraphXml
<w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $prng.InsertXML($pkg)
}

function Insert-ParagraphBefore($anchorText, $newParagraphXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "anchor text not found: $anchorText"
    }
    $p = $rng.Paragraphs(1)
    $prng = $p.Range
    $insPoint = $d.Range($prng.Start, $prng.Start)
    $pkg = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$newParagraphXml
<w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $insPoint.InsertXML($pkg)
}

# --- Change 1: TIENEN -> TIENE ---------------------------------------------
$para1 = @"
<w:p w14:paraId="3B264B53" w14:textId="6AE358F4" w:rsidR="0031725D" w:rsidRDefault="00000000" w:rsidP="00DC6A9D"><w:pPr><w:tabs><w:tab w:val="left" w:pos="142"/></w:tabs><w:spacing w:line="276" w:lineRule="auto"/><w:ind w:right="-660"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:b/></w:rPr><w:t xml:space="preserve">PRIMERA. </w:t></w:r><w:bookmarkStart w:id="6" w:name="_Hlk207403546"/><w:r w:rsidR="007D5BF2"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:b/></w:rPr><w:t>“</w:t></w:r><w:r w:rsidR="0031725D" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>{{SEXO_1}} PROMITENTE {{SEXO_2}}</w:t></w:r><w:r w:rsidR="007D5BF2"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>”</w:t></w:r><w:r w:rsidR="0031725D" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>, SE COMPROMETE A VENDER {{SEXO_</w:t></w:r><w:r w:rsidR="0031725D"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>12</w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="0031725D" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>}}PROMITENTE</w:t></w:r><w:r w:rsidR="00764139"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>S</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="0031725D" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t xml:space="preserve"> {{SEXO_</w:t></w:r><w:r w:rsidR="0031725D"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>10</w:t></w:r><w:r w:rsidR="0031725D" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>}}</w:t></w:r><w:r w:rsidR="007D5BF2"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>”</w:t></w:r><w:r w:rsidR="0031725D" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>, EL LOTE #{{IDENTIFICADOR_LOTE}} ({{LETRA_IDENTIFICADOR}}) UBICADO EN {{DIRECCION_PROYECTO_LOTE}} EL CUAL TIENE LAS SIGUIENTES MEDIDA</w:t></w:r><w:r w:rsidR="007C5C52"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>S</w:t></w:r><w:r w:rsidR="0031725D" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t xml:space="preserve"> Y COLINDANCIAS</w:t></w:r><w:bookmarkEnd w:id="6"/><w:r w:rsidR="0031725D"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>.</w:t></w:r></w:p>
"@
Set-ParagraphXml "SE COMPROMETE A VENDER" $para1

# --- Change 2: insert blank Prrafodelista paragraph -------------------------
$newEmptyPara = @"
<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:right="-660"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr></w:pPr></w:p>
"@
Insert-ParagraphBefore "EL PAGO DEBERÁ REALIZARSE DE LUNES A SÁBADO" $newEmptyPara

# --- Change 3: OBLIGADOS -> OBLIGAD{{SEXO_11}}S -----------------------------
$para3 = @"
<w:p w14:paraId="6D0850EB" w14:textId="309527BD" w:rsidR="00B65D29" w:rsidRDefault="00000000"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:tabs><w:tab w:val="left" w:pos="142"/></w:tabs><w:spacing w:after="0" w:line="276" w:lineRule="auto"/><w:ind w:left="426" w:right="-660"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr></w:pPr><w:bookmarkStart w:id="53" w:name="_Hlk206455471"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t>EFECTUARÁ</w:t></w:r><w:r w:rsidR="00DC6A9D"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t>N</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="003F2DFD" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">LO ESTABLECIDO SOBRE LA PENA CONVENCIONAL Y LO INDICADO EN EL INCISO D) DE LA NOVENA CLÁUSULA. SI EL MONTO ABONADO POR </w:t></w:r><w:r w:rsidR="007D5BF2"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t>“</w:t></w:r><w:r w:rsidR="001473D2"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>{{SEXO_9}}</w:t></w:r><w:r w:rsidR="003F2DFD" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> PROMITENTE</w:t></w:r><w:r w:rsidR="001473D2"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t>S</w:t></w:r><w:r w:rsidR="003F2DFD" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="001473D2"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>{{SEXO_10}}</w:t></w:r><w:r w:rsidR="007D5BF2"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>”</w:t></w:r><w:r w:rsidR="003F2DFD" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> SUPERA LA PENA CONVENCIONAL, </w:t></w:r><w:r w:rsidR="007D5BF2"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t>“</w:t></w:r><w:r w:rsidR="003F2DFD" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>{{SEXO_1}}</w:t></w:r><w:r w:rsidR="003F2DFD" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> PROMITENTE </w:t></w:r><w:r w:rsidR="003F2DFD" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>{{SEXO_2}}</w:t></w:r><w:r w:rsidR="007D5BF2"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>”</w:t></w:r><w:r w:rsidR="003F2DFD" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> LE DEVOLVERÁ LA DIFERENCIA; DE LO CONTRARIO, </w:t></w:r><w:r w:rsidR="007D5BF2"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t>“</w:t></w:r><w:r w:rsidR="001473D2"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>{{SEXO_9}}</w:t></w:r><w:r w:rsidR="003F2DFD" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> PROMITENTE</w:t></w:r><w:r w:rsidR="001473D2"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t>S</w:t></w:r><w:r w:rsidR="003F2DFD" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="001473D2"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>{{SEXO_10}}</w:t></w:r><w:r w:rsidR="007D5BF2"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/></w:rPr><w:t>”</w:t></w:r><w:r w:rsidR="003F2DFD" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> QUEDARÁ</w:t></w:r><w:r w:rsidR="006819B1"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t>N</w:t></w:r><w:r w:rsidR="003F2DFD" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> OBLIGAD</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t>{{SEXO_11</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t>}}</w:t></w:r><w:r w:rsidR="006819B1"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t>S</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="003F2DFD" w:rsidRPr="005F206A"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> A CUBRIR LA PENA CONVENCIONAL</w:t></w:r><w:bookmarkEnd w:id="53"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Arial Narrow"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r></w:p>
"@
Set-ParagraphXml "A CUBRIR LA PENA CONVENCIONAL" $para3

Write-Output "done"
